$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BugFix: Specific ions file accurate masses
$ws.Range("C14").Value = 184.0739
$ws.Range("C20").Value = 142.0269
$ws.Range("C31").Value = 186.0168

# Temporary fix: align the numeric format of C31 with the other corrected,
# right-aligned numeric cells in the column (matches style used by C20 etc.)
$ws.Range("C31").NumberFormat = "0.0000"
$ws.Range("C31").HorizontalAlignment = -4152 # xlRight
$ws.Range("C31").VerticalAlignment = -4108   # xlCenter

# Adjust the visible window / selection to reflect the author's saved view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select()
